$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: old "@EffectIcon" (Damage.ai/Protection.ai) -> new "EffectText" values ---
$ws.Range("C2").Value = "-1 HP"
$ws.Range("C3").Value = "-2 HP"
$ws.Range("C4").Value = "+1 HP"
$ws.Range("C5").Value = "+1 HP"
$ws.Range("C6").Value = "+2 HP"
$ws.Range("C7").Value = "-3 HP"
$ws.Range("C8").Value = "+3 HP"

# Values starting with +/- get the "quote-prefixed" style, matching the
# style already used for this kind of value elsewhere in the sheet (E3).
$ws.Range("E3").Copy()
$ws.Range("C2:C8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column D: old "EffectText" (x1/x2/x3) -> new "ComboText" values ---
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = "1 extra damage with a Poisoner"
$ws.Range("D4").Value = "Destroys one Poisoner"
$ws.Range("D5").Value = "Destroys one Saboteur"
$ws.Range("D6").Value = "Destroys one Assassin"
$ws.Range("D7").Value = "1 extra damage with a Poisoner"
$ws.Range("D8").Value = "Destroys one Assassin"

# --- Remove the old "@ComboIcon" / "ComboText" / "ComboPrefix" columns (E:G) ---
$ws.Range("E1:G8").EntireColumn.Delete()

# --- Header row: rename C1/D1 ---
$ws.Range("C1").Value = "EffectText"
$ws.Range("D1").Value = "ComboText"

# C1 is no longer an "@"-prefixed field, so it should use the plain header
# style (matching A1/D1) rather than the quote-prefixed header style.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column D width ---
$ws.Columns.Item(4).ColumnWidth = 37.3

# --- Selection ---
$ws.Range("D2").Select()
